# Automatic update of files.
#
# Rows 21, 22 and 24 of the sheet get their species-record data rotated
# (row21 -> row22, row22 -> row24, row24 -> row21), each record keeping
# its own "Taxonsorteringsordning" (column B) bumped by 14. Row 23 only
# gets its column B value bumped as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: becomes the "Svartvit taggsvamp / Phellodon connatus" record ---
$ws.Range("A21").Value = 111895157
$ws.Range("B21").Value = 90857
$ws.Range("E21").Value = 5448
$ws.Range("F21").Value = "Svartvit taggsvamp"
$ws.Range("G21").Value = "Phellodon connatus"
$ws.Range("H21").Value = "(Schultz) nom.prov"
$ws.Range("Q21").Value = 383311
$ws.Range("R21").Value = 6664460
# S21 stays 10 (unchanged)

# --- Row 22: becomes the "Motaggsvamp / Sarcodon squamosus" record ---
$ws.Range("A22").Value = 111895200
$ws.Range("B22").Value = 90837
$ws.Range("E22").Value = 5966
$ws.Range("F22").Value = "Motaggsvamp"
$ws.Range("G22").Value = "Sarcodon squamosus"
$ws.Range("H22").Value = "(Schaeff.) Quél."
$ws.Range("I22").ClearContents()
$ws.Range("L22").ClearContents()
$ws.Range("M22").ClearContents()
$ws.Range("Q22").Value = 383318
$ws.Range("R22").Value = 6664423
$ws.Range("S22").Value = 10

# --- Row 23: only the Taxonsorteringsordning (B) changes ---
$ws.Range("B23").Value = 90830

# --- Row 24: becomes the "Spillkråka / Dryocopus martius" record ---
$ws.Range("A24").Value = 111895144
$ws.Range("B24").Value = 56446
$ws.Range("E24").Value = 100049
$ws.Range("F24").Value = "Spillkråka"
$ws.Range("G24").Value = "Dryocopus martius"
$ws.Range("H24").Value = "(Linnaeus, 1758)"
$ws.Range("I24").Value = "'1"
$ws.Range("J24").ClearContents()
$ws.Range("M24").Value = "lockläte, övriga läten"
$ws.Range("Q24").Value = 383215
$ws.Range("R24").Value = 6664539
$ws.Range("S24").Value = 25
$ws.Range("AF24").ClearContents()
